$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.272327238179451
$ws.Range("C2").Value = 0.3048912486333797
$ws.Range("D2").Value = 0.7210945179870265
$ws.Range("E2").Value = 0.5333859586016987
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 4.831698963401555

$ws.Range("B3").Value = 3.272327238179451
$ws.Range("C3").Value = 0.3048912486333797
$ws.Range("D3").Value = 0.7210945179870265
$ws.Range("E3").Value = 0.5333859586016987
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 4.831698963401555

$ws.Range("B4").Value = 3.272327238179451
$ws.Range("C4").Value = 1.626987699542094
$ws.Range("D4").Value = 0.1496068669990043
$ws.Range("E4").Value = 0.5333859586016987
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 5.582307763322248

$ws.Range("B5").Value = 3.272327238179451
$ws.Range("C5").Value = 1.626987699542094
$ws.Range("D5").Value = 0.7210945179870265
$ws.Range("E5").Value = 0.5333859586016987
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 6.15379541431027

$ws.Range("B6").Value = 3.272327238179451
$ws.Range("C6").Value = 1.626987699542094
$ws.Range("D6").Value = 3.223369029078222
$ws.Range("E6").Value = 0.5333859586016987
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 8.656069925401464
